$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): relabel / reorder columns F:L, drop column M ---
$ws.Range("F1").Value = "Maintainability"
$ws.Range("G1").Value = "Reliability"
$ws.Range("H1").Value = "Total_QR"
$ws.Range("I1").Value = "LOC"
$ws.Range("J1").Value = "Repository"
$ws.Range("K1").Value = "Version"
$ws.Range("L1").Value = "Data"

# --- Data rows 2:10 ---
# Columns F:L are rebuilt from the previous F:M values (maintainability,
# reliability, total, ncloc, repository, version, version-again) and the
# old column M is dropped entirely.
$data = @(
    @{ Row = 2;  F = 0.396;              G = 0.5;                H = 0.896;              I = 200;  J = "Archives"; K = "03-10-2021-15-47" },
    @{ Row = 3;  F = 0.33;               G = 0.5;                H = 0.83;               I = 428;  J = "Archives"; K = "03-11-2021-02-10" },
    @{ Row = 4;  F = 0.396;              G = 0.5;                H = 0.896;              I = 109;  J = "Archives"; K = "13-09-2021-15-00" },
    @{ Row = 5;  F = 0.396;              G = 0.5;                H = 0.896;              I = 112;  J = "Profile";  K = "03-10-2021-15-48" },
    @{ Row = 6;  F = 0.396;              G = 0.5;                H = 0.896;              I = 112;  J = "Profile";  K = "13-09-2021-14-00" },
    @{ Row = 7;  F = 0.3225;             G = 0.5;                H = 0.8225;             I = 1401; J = "Frontend"; K = "03-10-2021-15-49" },
    @{ Row = 8;  F = 0.3131632653061225; G = 0.4785714285714285; H = 0.791734693877551;  I = 4888; J = "Frontend"; K = "03-11-2021-02-08" },
    @{ Row = 9;  F = 0.33;               G = 0.5;                H = 0.8300000000000001; I = 35;   J = "Frontend"; K = "13-09-2021-20-00" },
    @{ Row = 10; F = 0.3131632653061225; G = 0.4857142857142857; H = 0.7988775510204081; I = 4858; J = "Frontend"; K = "17-10-2021-15-30" }
)

foreach ($r in $data) {
    $ws.Cells.Item($r.Row, 6).Value  = $r.F    # F - Maintainability
    $ws.Cells.Item($r.Row, 7).Value  = $r.G    # G - Reliability
    $ws.Cells.Item($r.Row, 8).Value  = $r.H    # H - Total_QR
    $ws.Cells.Item($r.Row, 9).Value  = $r.I    # I - LOC
    $ws.Cells.Item($r.Row, 10).Value = $r.J    # J - Repository
    $ws.Cells.Item($r.Row, 11).Value = $r.K    # K - Version
    $ws.Cells.Item($r.Row, 12).Value = $r.K    # L - Data (duplicate of Version)
}

# Old column M (total) is no longer part of the table; remove it entirely.
$ws.Columns.Item(13).Delete()
